# Update "想去人数" (column F) figures on each worksheet to match the
# newly generated output (commit: "Update gh-pages to output generated at a56beed").

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 4-29
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 747
$ws1.Range("F5").Value = 491
$ws1.Range("F6").Value = 251
$ws1.Range("F7").Value = 8495
$ws1.Range("F8").Value = 224
$ws1.Range("F9").Value = 603
$ws1.Range("F10").Value = 1412
$ws1.Range("F11").Value = 21
$ws1.Range("F12").Value = 12
$ws1.Range("F13").Value = 2056
$ws1.Range("F15").Value = 3397
$ws1.Range("F16").Value = 193
$ws1.Range("F17").Value = 53
$ws1.Range("F18").Value = 75
$ws1.Range("F19").Value = 139
$ws1.Range("F20").Value = 204
$ws1.Range("F21").Value = 157
$ws1.Range("F22").Value = 2
$ws1.Range("F23").Value = 70
$ws1.Range("F24").Value = 221
$ws1.Range("F25").Value = 92
$ws1.Range("F26").Value = 1014
$ws1.Range("F27").Value = 352
$ws1.Range("F28").Value = 4213
$ws1.Range("F29").Value = 19

# Sheet "演出": row 3
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 23

# Sheet "本地生活": rows 2-3
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 156
$ws3.Range("F3").Value = 820

# Sheet "全部类型": rows 2-33
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 156
$ws4.Range("F4").Value = 820
$ws4.Range("F6").Value = 23
$ws4.Range("F8").Value = 747
$ws4.Range("F9").Value = 491
$ws4.Range("F10").Value = 251
$ws4.Range("F11").Value = 8495
$ws4.Range("F12").Value = 224
$ws4.Range("F13").Value = 603
$ws4.Range("F14").Value = 1412
$ws4.Range("F15").Value = 21
$ws4.Range("F16").Value = 12
$ws4.Range("F17").Value = 2056
$ws4.Range("F19").Value = 3397
$ws4.Range("F20").Value = 193
$ws4.Range("F21").Value = 53
$ws4.Range("F22").Value = 75
$ws4.Range("F23").Value = 139
$ws4.Range("F24").Value = 204
$ws4.Range("F25").Value = 157
$ws4.Range("F26").Value = 2
$ws4.Range("F27").Value = 70
$ws4.Range("F28").Value = 221
$ws4.Range("F29").Value = 92
$ws4.Range("F30").Value = 1014
$ws4.Range("F31").Value = 352
$ws4.Range("F32").Value = 4213
$ws4.Range("F33").Value = 19

$wb.Save()
